# Auto-generated Excel COM-interop script
# Updates crypto price ("D" column) and 1h volume change ("E" column)
# values on Sheet1 to match the refreshed data snapshot from the
# GitHub Actions run. Values are plain text (numeric-looking strings),
# matching the original workbook's storage format, so each assignment
# uses a leading apostrophe to force text entry (prevents Excel from
# auto-converting "304.64" / "0.83%" into a Number/Percentage value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.64"
$ws.Range("E2").Value = "'0.83%"
$ws.Range("D3").Value = "'35.93"
$ws.Range("E3").Value = "'2.02%"
$ws.Range("D4").Value = "'5.107"
$ws.Range("E4").Value = "'1.27%"
$ws.Range("D5").Value = "'0.08077"
$ws.Range("E5").Value = "'1.17%"
$ws.Range("D6").Value = "'1.935"
$ws.Range("E6").Value = "'-0.05%"
$ws.Range("D7").Value = "'4.186"
$ws.Range("E7").Value = "'3.44%"
$ws.Range("D8").Value = "'7.751"
$ws.Range("E8").Value = "'0.14%"
$ws.Range("D9").Value = "'0.9270"
$ws.Range("E9").Value = "'0.43%"
$ws.Range("D10").Value = "'0.1356"
$ws.Range("E10").Value = "'4.64%"
$ws.Range("D11").Value = "'0.1904"
$ws.Range("E11").Value = "'3.15%"
$ws.Range("D12").Value = "'0.09179"
$ws.Range("E12").Value = "'-5.09%"
$ws.Range("D13").Value = "'0.03416"
$ws.Range("E13").Value = "'-5.99%"
$ws.Range("D14").Value = "'0.09824"
$ws.Range("E14").Value = "'-0.28%"
$ws.Range("D15").Value = "'0.001414"
$ws.Range("E15").Value = "'1.40%"
$ws.Range("D16").Value = "'0.005770"
$ws.Range("E16").Value = "'-0.87%"
$ws.Range("D17").Value = "'3.555"
$ws.Range("E17").Value = "'1.44%"
$ws.Range("D18").Value = "'2.967"
$ws.Range("E18").Value = "'1.91%"
$ws.Range("D19").Value = "'0.3453"
$ws.Range("E19").Value = "'0.68%"
$ws.Range("D20").Value = "'0.1332"
$ws.Range("E20").Value = "'1.65%"
$ws.Range("D21").Value = "'4.906"
$ws.Range("E21").Value = "'-2.78%"
$ws.Range("D22").Value = "'0.2602"
$ws.Range("E22").Value = "'8.34%"
$ws.Range("D23").Value = "'0.04399"
$ws.Range("E23").Value = "'-2.78%"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'0.31%"
$ws.Range("D25").Value = "'0.004808"
$ws.Range("E25").Value = "'-0.15%"
$ws.Range("E26").Value = "'4.17%"
$ws.Range("D27").Value = "'0.0003133"
$ws.Range("E27").Value = "'4.23%"
$ws.Range("E39").Value = "'5.69%"
$ws.Range("D40").Value = "'0.04904"
$ws.Range("E40").Value = "'4.34%"
$ws.Range("D41").Value = "'0.007621"
$ws.Range("E41").Value = "'1.26%"
$ws.Range("D42").Value = "'0.01026"
$ws.Range("E42").Value = "'5.86%"
$ws.Range("E43").Value = "'4.20%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-0.31%"
$ws.Range("D45").Value = "'0.01089"
$ws.Range("E45").Value = "'0.58%"
$ws.Range("D46").Value = "'0.00006400"
$ws.Range("E46").Value = "'2.67%"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.41%"
$ws.Range("D49").Value = "'0.001192"
$ws.Range("E49").Value = "'-19.97%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.03%"
